$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4
$ws.Range("A2").Value = 22
$ws.Range("B2").Value = 3

$ws.Range("A3").Value = 21
$ws.Range("B3").Value = 2

$ws.Range("A4").Value = 12
$ws.Range("B4").Value = 1

# Add new row 5, copying style from A4/B4's previous style (bold, border) for A5
$ws.Range("A5").Value = 11
$ws.Range("B5").Value = 1

# Apply style to A5 matching A2:A4 (copy format from A4)
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122) # xlPasteFormats
